{"js": "// The document's auto-managed \"_GoBack\" bookmark (Word's \"last edit\n// location\" marker) moves from the middle of the \"Buka command prompt\n// with ruby on rails...\" paragraph to the end of the final paragraph,\n// and that final paragraph's single run is split in two at the point\n// the edit happened (\"...terdapat pada |web, silahkan login...\"),\n// matching a user having typed/edited right before \"web, silahkan...\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Remove the old \"_GoBack\" bookmark (it sat between \"with \" and\n//    \"ruby on rails ...\" in the \"Buka command prompt\" step).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the final paragraph (\"Untuk mengeksplorasi ...\") and split\n//    its run into two runs at \"web, silahkan login ke dalam web\".\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Untuk mengeksplorasi dan menggunakan fitur-fitur\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const searchResults = targetParagraph.search(\"web, silahkan login ke dalam web\", { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length > 0) {\n    // Force the run boundary to appear exactly at the split point by\n    // inserting (then immediately deleting) a scratch bookmark there.\n    const splitPoint = searchResults.items[0].getRange(\"Start\");\n    splitPoint.insertBookmark(\"_TempSplitMarker\");\n    await context.sync();\n\n    context.document.deleteBookmark(\"_TempSplitMarker\");\n    await context.sync();\n  }\n\n  // 3) Re-insert \"_GoBack\" at the very end of that paragraph.\n  const endRange = targetParagraph.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The document's auto-managed \"_GoBack\" bookmark (Word's \"last edit\n# location\" marker) moves from the middle of the \"Buka command prompt\n# with ruby on rails...\" paragraph to the end of the final paragraph,\n# and that final paragraph's single run is split in two at the point\n# the edit happened (\"...terdapat pada |web, silahkan login...\"),\n# matching a user having typed/edited right before \"web, silahkan...\".\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark (it sat between \"with \" and\n#    \"ruby on rails ...\" in the \"Buka command prompt\" step).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Locate the final paragraph (\"Untuk mengeksplorasi ...\") and split\n#    its run into two runs at \"web, silahkan login ke dalam web\".\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Untuk mengeksplorasi dan menggunakan fitur-fitur*\") {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -ne $null) {\n    $pRange = $targetParagraph.Range\n\n    # 2a. Force a run boundary right before \"web, silahkan login ke dalam\n    #     web\" using a scratch bookmark that is immediately deleted again\n    #     (the split itself is what survives, not the bookmark).\n    $found = $pRange.Duplicate\n    $ok = $found.Find.Execute(\"web, silahkan login ke dalam web\")\n    if ($ok) {\n        $splitPoint = $found.Duplicate\n        $splitPoint.Collapse(1)   # wdCollapseStart\n        $d.Bookmarks.Add(\"_TempSplitMarker\", $splitPoint)\n        $d.Bookmarks.Item(\"_TempSplitMarker\").Delete()\n    }\n\n    # 2b. Insert \"_GoBack\" at the very end of that paragraph's text\n    #     (right before its paragraph mark). A collapsed range placed\n    #     exactly at that boundary does not anchor reliably, so a\n    #     throwaway character is inserted, the bookmark is wrapped\n    #     around it, and then the character is cleared through the\n    #     bookmark's own range -- leaving a clean, correctly anchored,\n    #     zero-length bookmark at the paragraph's end.\n    $endRange = $targetParagraph.Range\n    $endRange.Collapse(0)      # wdCollapseEnd\n    $endRange.MoveEnd(1, -1) | Out-Null\n    $endRange.InsertAfter(\"X\")\n    $d.Bookmarks.Add(\"_GoBack\", $endRange)\n    $d.Bookmarks.Item(\"_GoBack\").Range.Text = \"\"\n}\n"}
